{"js": "// Change: \"...extra information regarding the datasets units, concepts, and\n// labeling choices is mentioned.\" -> \"...extra information regarding the\n// dataset's units, concepts, and labeling choices are mentioned.\"\n//\n// Two independent word-level fixes inside the same sentence:\n//   1. \"datasets\" -> \"dataset's\"   (singular possessive)\n//   2. \"is mentioned\" -> \"are mentioned\"  (subject/verb agreement with the\n//      now-plural-sounding compound subject \"units, concepts, and labeling\n//      choices\")\n\nconst body = context.document.body;\n\n// Target the unique sentence fragment so the generic word \"datasets\" that\n// appears elsewhere in the document (e.g. \"a number of datasets and\n// processing steps\") is left untouched.\nconst phrase =\n  \"regarding the datasets units, concepts, and labeling choices is mentioned\";\nconst phraseResults = body.search(phrase, { matchCase: true, matchWholeWord: false });\nphraseResults.load(\"items\");\nawait context.sync();\n\nif (phraseResults.items.length === 0) {\n  throw new Error(\"Could not find the target sentence fragment to replace.\");\n}\n\n// 1) \"datasets\" -> \"dataset\\u2019s\", scoped to the matched sentence only.\nconst datasetResults = phraseResults.items[0].search(\"datasets\", {\n  matchCase: true,\n  matchWholeWord: true\n});\ndatasetResults.load(\"items\");\nawait context.sync();\n\nif (datasetResults.items.length === 0) {\n  throw new Error('Could not find \"datasets\" to replace.');\n}\n\ndatasetResults.items[0].insertText(\"dataset\\u2019s\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) \" is mentioned\" -> \" are mentioned\". The previous range object is now\n// stale (its backing text changed), so re-locate the updated sentence\n// fragment from the body before scoping the second search.\nconst updatedPhrase =\n  \"regarding the dataset\\u2019s units, concepts, and labeling choices is mentioned\";\nconst updatedPhraseResults = body.search(updatedPhrase, {\n  matchCase: true,\n  matchWholeWord: false\n});\nupdatedPhraseResults.load(\"items\");\nawait context.sync();\n\nif (updatedPhraseResults.items.length === 0) {\n  throw new Error(\"Could not re-find the updated sentence fragment.\");\n}\n\nconst isMentionedResults = updatedPhraseResults.items[0].search(\"is mentioned\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nisMentionedResults.load(\"items\");\nawait context.sync();\n\nif (isMentionedResults.items.length === 0) {\n  throw new Error('Could not find \"is mentioned\" to replace.');\n}\n\nisMentionedResults.items[0].insertText(\"are mentioned\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Change: \"...extra information regarding the datasets units, concepts, and\n# labeling choices is mentioned.\" -> \"...extra information regarding the\n# dataset's units, concepts, and labeling choices are mentioned.\"\n#\n# Two independent word-level fixes inside the same sentence:\n#   1. \"datasets\" -> \"dataset's\"   (singular possessive)\n#   2. \"is mentioned\" -> \"are mentioned\"  (subject/verb agreement)\n#\n# The word \"datasets\" also appears elsewhere in the document (e.g. \"a number\n# of datasets and processing steps\"), so the fix is scoped to the one\n# paragraph that contains the target sentence instead of running Find over\n# the whole document.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*extra information regarding the datasets units, concepts, and labeling choices is mentioned*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the target paragraph.\"\n}\n\n# 1) \"datasets\" -> \"dataset's\" (curly/typographic apostrophe U+2019,\n# matching the rest of the document's punctuation style), scoped to this\n# paragraph only.\n$rng1 = $target.Range\n$found1 = $rng1.Find.Execute(\"datasets\", $true, $true, $false, $false, $false, $true, 1, $false, \"dataset\u2019s\", 2)\nif (-not $found1) {\n    throw \"Could not find 'datasets' to replace.\"\n}\n\n# 2) \"is mentioned\" -> \"are mentioned\", scoped to the same paragraph.\n$rng2 = $target.Range\n$found2 = $rng2.Find.Execute(\"is mentioned\", $true, $false, $false, $false, $false, $true, 1, $false, \"are mentioned\", 2)\nif (-not $found2) {\n    throw \"Could not find 'is mentioned' to replace.\"\n}\n"}
